# Generate Report for Handback
#
# The "97e04619-991d-46e2-8ce3-16b4a0b6d44a" file has now been handed back
# (status flips from "Ready for handoff" to "Handed back: in sync with
# en-US") while "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a" is still only ready
# for handoff. Because the report is regenerated/re-sorted, the two data
# rows on every sheet swap places: row 2 now belongs to 97e04619 and row 3
# to 03ca7e5d. The per-language sheets additionally grow two new columns
# of data (F: Latest Target File, G: Latest Handback File) plus an H
# (Latest Handback DateTime) value for the file that was just handed back.

$wb = $excel.ActiveWorkbook

$mdA = "97e04619-991d-46e2-8ce3-16b4a0b6d44a.md"
$mdB = "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.md"

$urlMdA = "https://github.com/OpenLocalizationTest/oltest/blob/766bebe7dff4ace04aa0ba942229366132bf4c29/e2e/97e04619-991d-46e2-8ce3-16b4a0b6d44a.md"
$urlMdB = "https://github.com/OpenLocalizationTest/oltest/blob/dccff3bde1564045b3cb46c0120653a4ec53d861/e2e/03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.md"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

# Row 2 -> 97e04619 (handed back)
$ws.Range("A2").Value = $mdA
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("C2").Value = $statusHandedBack
$ws.Range("D2").Value = "2016-26-18 12:26:27"

# Row 3 -> 03ca7e5d (still only ready for handoff)
$ws.Range("A3").Value = $mdB
$ws.Range("B3").Value = $statusReady
$ws.Range("C3").Value = $statusReady
$ws.Range("D3").Value = "2016-26-18 12:26:04"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, $mdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMdB, [System.Type]::Missing, [System.Type]::Missing, $mdB) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$xlfHandoffA_zh = "97e04619-991d-46e2-8ce3-16b4a0b6d44a.54537a46c86a074c049c2dfc438b8658616df32c.zh-cn.xlf"
$xlfHandoffB_zh = "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.81200effe507a49bc7034878fb2a2f18ca8e9f06.zh-cn.xlf"

$urlXlfHandoffA_zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2b1e02b161b8430809e1fbd38102cbd65c2be66/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97e04619-991d-46e2-8ce3-16b4a0b6d44a.54537a46c86a074c049c2dfc438b8658616df32c.zh-cn.xlf"
$urlXlfHandoffB_zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aed93a2465ab09b6dea504422ca354ab41d8d08/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.81200effe507a49bc7034878fb2a2f18ca8e9f06.zh-cn.xlf"
$urlXlfHandbackA_zh = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/97e04619991d46e28ce316b4a0b6d44a54537a4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/97e04619-991d-46e2-8ce3-16b4a0b6d44a.54537a46c86a074c049c2dfc438b8658616df32c.zh-cn.xlf"

# Row 2 -> 97e04619 (handed back: handoff + target + handback all present)
$ws.Range("A2").Value = $mdA
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $statusHandedBack
$ws.Range("D2").Value = $xlfHandoffA_zh
$ws.Range("E2").Value = "2016-03-18 12:26:23"
$ws.Range("F2").Value = $mdA
$ws.Range("G2").Value = $xlfHandoffA_zh
$ws.Range("H2").Value = "2016-03-18 12:26:47"
$ws.Range("I2").Value = "Include"

# Row 3 -> 03ca7e5d (only handed off so far, no handback yet)
$ws.Range("A3").Value = $mdB
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $statusReady
$ws.Range("D3").Value = $xlfHandoffB_zh
$ws.Range("E3").Value = "2016-03-18 12:26:00"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

Style-AsHyperlink $ws.Range("F2")
Style-AsHyperlink $ws.Range("G2")

$ws.Hyperlinks.Add($ws.Range("A2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, $mdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $urlXlfHandoffA_zh, [System.Type]::Missing, [System.Type]::Missing, $xlfHandoffA_zh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, $mdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), $urlXlfHandbackA_zh, [System.Type]::Missing, [System.Type]::Missing, $xlfHandoffA_zh) | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), $urlMdB, [System.Type]::Missing, [System.Type]::Missing, $mdB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $urlMdB, [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $urlXlfHandoffB_zh, [System.Type]::Missing, [System.Type]::Missing, $xlfHandoffB_zh) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$xlfHandoffA_de = "97e04619-991d-46e2-8ce3-16b4a0b6d44a.54537a46c86a074c049c2dfc438b8658616df32c.de-de.xlf"
$xlfHandoffB_de = "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.81200effe507a49bc7034878fb2a2f18ca8e9f06.de-de.xlf"

$urlXlfHandoffA_de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c8d0005c833f45ee4fd12c1cdae60468f9af2c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97e04619-991d-46e2-8ce3-16b4a0b6d44a.54537a46c86a074c049c2dfc438b8658616df32c.de-de.xlf"
$urlXlfHandoffB_de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4cb153bd301cd3c3491a6674b162bb1c8b1d3ff/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.81200effe507a49bc7034878fb2a2f18ca8e9f06.de-de.xlf"
$urlXlfHandbackA_de = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/97e04619991d46e28ce316b4a0b6d44a54537a4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/97e04619-991d-46e2-8ce3-16b4a0b6d44a.54537a46c86a074c049c2dfc438b8658616df32c.de-de.xlf"

# Row 2 -> 97e04619 (handed back: handoff + target + handback all present)
$ws.Range("A2").Value = $mdA
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $statusHandedBack
$ws.Range("D2").Value = $xlfHandoffA_de
$ws.Range("E2").Value = "2016-03-18 12:26:27"
$ws.Range("F2").Value = $mdA
$ws.Range("G2").Value = $xlfHandoffA_de
$ws.Range("H2").Value = "2016-03-18 12:26:52"
$ws.Range("I2").Value = "Include"

# Row 3 -> 03ca7e5d (only handed off so far, no handback yet)
$ws.Range("A3").Value = $mdB
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $statusReady
$ws.Range("D3").Value = $xlfHandoffB_de
$ws.Range("E3").Value = "2016-03-18 12:26:04"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

Style-AsHyperlink $ws.Range("F2")
Style-AsHyperlink $ws.Range("G2")

$ws.Hyperlinks.Add($ws.Range("A2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, $mdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $urlXlfHandoffA_de, [System.Type]::Missing, [System.Type]::Missing, $xlfHandoffA_de) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), $urlMdA, [System.Type]::Missing, [System.Type]::Missing, $mdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), $urlXlfHandbackA_de, [System.Type]::Missing, [System.Type]::Missing, $xlfHandoffA_de) | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), $urlMdB, [System.Type]::Missing, [System.Type]::Missing, $mdB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $urlMdB, [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $urlXlfHandoffB_de, [System.Type]::Missing, [System.Type]::Missing, $xlfHandoffB_de) | Out-Null

Write-Output "Report regenerated for handback."
